$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

# --- Row 2 ---
Set-TextValue "A2" "1"
Set-TextValue "B2" "00100102010000009206"
Set-TextValue "C2" "3694000010100"
$ws.Range("D2").Value = "Tapa #12-14-16-21 Plast BF 40x25 unid Dynapack"
# E2 (ImpuestoPorcentaje) unchanged: 13
$ws.Range("F2").Value = "DESCUENTO"
# G2, H2 unchanged
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 10391.33
$ws.Range("K2").Value = 1558.6995
$ws.Range("L2").Value = 1148.24197
$ws.Range("M2").Value = 8832.630499999999
$ws.Range("N2").Value = 9980.87247
Set-TextValue "O2" "3694046"
$ws.Range("P2").Value = "Unid"
$ws.Range("Q2").ClearContents()
$ws.Range("R2").Value = 10391.33
$ws.Range("S2").Value = "factura"

# --- Row 3 ---
Set-TextValue "A3" "2"
Set-TextValue "B3" "00100102010000009206"
Set-TextValue "C3" "3219905000200"
$ws.Range("D3").Value = "Vasos #22 EB Carton BF 40x25 unid"
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = "DESCUENTO"
# G3 unchanged
$ws.Range("H3").Value = 8
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 33408
$ws.Range("K3").Value = 5011.2
$ws.Range("L3").Value = 3691.584
$ws.Range("M3").Value = 28396.8
$ws.Range("N3").Value = 32088.384
Set-TextValue "O3" "3219933"
$ws.Range("P3").Value = "Unid"
$ws.Range("Q3").ClearContents()
$ws.Range("R3").Value = 33408
$ws.Range("S3").Value = "factura"

$wb.Save()
